$d = $word.ActiveDocument

# --- 1. Add new paragraph style "AbstractTitle" ---------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. Abstract style: spacing before 300 -> 100 --------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. ImportTok character style gains green bold colour ------------------
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 32768
$importTok.Font.Bold = $true

# --- 4. BuiltInTok character style gains green colour -----------------------
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768

Write-Output "done"
